$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (row index is 1-based)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "375"
$t.Cell(6,1).Range.Text  = "0.00083"
$t.Cell(7,1).Range.Text  = "0.00018"
$t.Cell(8,1).Range.Text  = "0.00005"
$t.Cell(9,1).Range.Text  = "0.00029"
$t.Cell(10,1).Range.Text = "0.00044"
$t.Cell(11,1).Range.Text = "0.00050"
$t.Cell(12,1).Range.Text = "0.07708"

# Collapse the tab-separated multi-run rows down to a single value each
$t.Cell(44,1).Range.Text = "99.97"
$t.Cell(45,1).Range.Text = "0.08"
$t.Cell(46,1).Range.Text = "234"
